$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1864.6666
$ws.Range("I19").Value = 2755.1
$ws.Range("J19").Value = 1228.6428
$ws.Range("K19").Value = 2755.1
$ws.Range("L19").Value = 1228.6428
$ws.Range("M19").Value = -2580.1
$ws.Range("N19").Value = -1578.6428

$ws.Range("H86").Value = 2944.6667
$ws.Range("I86").Value = 2700
$ws.Range("J86").Value = 3434
$ws.Range("K86").Value = 2700
$ws.Range("L86").Value = 3434
$ws.Range("M86").Value = -1577
$ws.Range("N86").Value = -5680

$ws.Range("H89").Value = 2944.6667
$ws.Range("I89").Value = 2700
$ws.Range("J89").Value = 3434
$ws.Range("K89").Value = 13500
$ws.Range("L89").Value = 17170
$ws.Range("M89").Value = -7884
$ws.Range("N89").Value = -28402

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""

$ws.Range("H137").Value = 3961.075
$ws.Range("I137").Value = 4532.552
$ws.Range("K137").Value = 13597.656
$ws.Range("M137").Value = -11047.656

$ws.Range("H141").Value = 6227.778
$ws.Range("I141").Value = 7720.7144
$ws.Range("J141").Value = 1002.5
$ws.Range("K141").Value = 23162.1432
$ws.Range("L141").Value = 3007.5
$ws.Range("M141").Value = -17982.1432
$ws.Range("N141").Value = -13367.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").Value = ""

$ws.Range("H32").Value = 13166.022
$ws.Range("I32").Value = 8077.758
$ws.Range("K32").Value = 8077.758
$ws.Range("M32").Value = -7790.758

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = ""

$ws.Range("H63").Value = 3229.5
$ws.Range("I63").Value = 1899.2858
$ws.Range("J63").Value = 6333.3335
$ws.Range("K63").Value = 1899.2858
$ws.Range("L63").Value = 6333.3335
$ws.Range("M63").Value = -1213.2858
$ws.Range("N63").Value = -7705.3335

$ws.Range("H66").Value = 3229.5
$ws.Range("I66").Value = 1899.2858
$ws.Range("J66").Value = 6333.3335
$ws.Range("K66").Value = 9496.429
$ws.Range("L66").Value = 31666.6675
$ws.Range("M66").Value = -6064.429
$ws.Range("N66").Value = -38530.6675

$ws.Range("H74").Value = 4767403.5
$ws.Range("I74").Value = 5882998.5
$ws.Range("J74").Value = 26125
$ws.Range("K74").Value = 5882998.5
$ws.Range("L74").Value = 26125
$ws.Range("M74").Value = -5882124.5
$ws.Range("N74").Value = -27873

$ws.Range("H77").Value = 4767403.5
$ws.Range("I77").Value = 5882998.5
$ws.Range("J77").Value = 26125
$ws.Range("K77").Value = 29414992.5
$ws.Range("L77").Value = 130625
$ws.Range("M77").Value = -29410624.5
$ws.Range("N77").Value = -139361

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = ""

$ws.Range("H122").Value = 2576.8462
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 18360
$ws.Range("J35").Value = 18360
$ws.Range("L35").Value = 18360
$ws.Range("N35").Value = -18980

$ws.Range("H82").Value = 13673.4375
$ws.Range("I82").Value = 9891
$ws.Range("J82").Value = 25020.75
$ws.Range("K82").Value = 9891
$ws.Range("L82").Value = 25020.75
$ws.Range("M82").Value = -9508
$ws.Range("N82").Value = -25786.75

$ws.Range("H85").Value = 13673.4375
$ws.Range("I85").Value = 9891
$ws.Range("J85").Value = 25020.75
$ws.Range("K85").Value = 9891
$ws.Range("L85").Value = 25020.75
$ws.Range("M85").Value = -8565
$ws.Range("N85").Value = -27672.75

$ws.Range("H94").Value = 1277.9565
$ws.Range("I94").Value = 1104.421
$ws.Range("K94").Value = 1104.421
$ws.Range("M94").Value = -653.421

$ws.Range("H132").Value = 60780
$ws.Range("J132").Value = 60780
$ws.Range("L132").Value = 60780
$ws.Range("N132").Value = -70900

$ws.Range("H140").Value = 77390
$ws.Range("J140").Value = 77390
$ws.Range("L140").Value = 77390
$ws.Range("N140").Value = -87750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5557544
$ws.Range("I31").Value = 1307.1904
$ws.Range("J31").Value = 18522096
$ws.Range("K31").Value = 1307.1904
$ws.Range("L31").Value = 18522096
$ws.Range("M31").Value = -1012.1904
$ws.Range("N31").Value = -18522686

$ws.Range("H34").Value = 5557544
$ws.Range("I34").Value = 1307.1904
$ws.Range("J34").Value = 18522096
$ws.Range("K34").Value = 1307.1904
$ws.Range("L34").Value = 18522096
$ws.Range("M34").Value = -1105.1904
$ws.Range("N34").Value = -18522500

$ws.Range("H70").Value = 35000
$ws.Range("J70").Value = 35000
$ws.Range("L70").Value = 35000
$ws.Range("N70").Value = -35630

$ws.Range("H73").Value = 35000
$ws.Range("J73").Value = 35000
$ws.Range("L73").Value = 35000
$ws.Range("N73").Value = -37184

$ws.Range("H97").Value = 13900
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 13900
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 13900
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = -15882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 80.23077000000001
$ws.Range("I40").Value = 64.333336
$ws.Range("J40").Value = 116
$ws.Range("K40").Value = 257.333344
$ws.Range("L40").Value = 464
$ws.Range("M40").Value = -188.333344
$ws.Range("N40").Value = -602

$ws.Range("H92").Value = 476.0909
$ws.Range("I92").Value = 490
$ws.Range("J92").Value = 464.5
$ws.Range("K92").Value = 1470
$ws.Range("L92").Value = 1393.5
$ws.Range("M92").Value = -222
$ws.Range("N92").Value = -3889.5

$ws.Range("H113").Value = 501.1111
$ws.Range("J113").Value = 461.13333
$ws.Range("L113").Value = 1383.39999
$ws.Range("N113").Value = -5723.39999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 122190.6
$ws.Range("I80").Value = 2180
$ws.Range("J80").Value = 242201.2
$ws.Range("K80").Value = 2180
$ws.Range("L80").Value = 242201.2
$ws.Range("M80").Value = -1182
$ws.Range("N80").Value = -244197.2

$ws.Range("H83").Value = 122190.6
$ws.Range("I83").Value = 2180
$ws.Range("J83").Value = 242201.2
$ws.Range("K83").Value = 10900
$ws.Range("L83").Value = 1211006
$ws.Range("M83").Value = -5908
$ws.Range("N83").Value = -1220990

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2972.32
$ws.Range("I40").Value = 2887.8333
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2887.8333
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2751.8333
$ws.Range("N40").Value = -5272

$ws.Range("H139").Value = 45665
$ws.Range("J139").Value = 45665
$ws.Range("L139").Value = 45665
$ws.Range("N139").Value = -55945

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1567.1
$ws.Range("I126").Value = 1818.8
$ws.Range("K126").Value = 5456.4
$ws.Range("M126").Value = -2986.4
